$wb = $excel.ActiveWorkbook

# --- Sheet "SQL": add rows 14 and 15 ---
$wsSql = $wb.Worksheets.Item("SQL")

$wsSql.Cells.Item(14, 1).Value = 1163
$wsSql.Cells.Item(14, 2).Value = "TestTestTestTestTestTestTest"
$wsSql.Cells.Item(14, 3).Value = "TestTestTestTestTestTestTest"
$wsSql.Cells.Item(14, 4).Value = "TestTestTestTestTestTestTest"

$wsSql.Cells.Item(15, 1).Value = 1164
$wsSql.Cells.Item(15, 2).Value = "TestTestTestTestTestTestTest"
$wsSql.Cells.Item(15, 3).Value = "TestTestTestTestTestTestTest"
$wsSql.Cells.Item(15, 4).Value = "TestTestTestTestTestTestTest"

# --- Sheet "Python": add rows 31 and 32 ---
$wsPy = $wb.Worksheets.Item("Python")

$wsPy.Cells.Item(31, 1).Value = 2091
$wsPy.Cells.Item(31, 2).Value = "TestTestTestTestTestTestTest"
$wsPy.Cells.Item(31, 3).Value = "TestTestTestTestTestTestTest"
$wsPy.Cells.Item(31, 4).Value = "TestTestTestTestTestTestTest"

$wsPy.Cells.Item(32, 1).Value = 2092
$wsPy.Cells.Item(32, 2).Value = "пвавапавпав"
$wsPy.Cells.Item(32, 3).Value = "павпвапавпвап"
$wsPy.Cells.Item(32, 4).Value = "вапввввввввв"

# --- Sheet "Links": update row 1 ---
$wsLinks = $wb.Worksheets.Item("Links")

$wsLinks.Cells.Item(1, 1).Value = 28
$wsLinks.Cells.Item(1, 2).Value = "TestTestTestTestTestTestTest"
$wsLinks.Cells.Item(1, 3).Value = "TestTestTestTestTestTestTest"
$wsLinks.Cells.Item(1, 4).Value = "TestTestTestTestTestTestTest"
